# Auto-generated Word COM-interop script
# Applies proofErr (spell/grammar check marks) splitting, a bookmark move,
# and new bold '${cloneValue}' template-block content, per the target diff.

$d = $word.ActiveDocument

function Set-ParaRuns {
    param($doc, $index, [string]$innerXml)
    $full = $doc.Paragraphs($index).Range
    $sub = $doc.Range($full.Start, $full.End - 1)
    $ns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
    $xml = "<w:p xmlns:w='$ns'>" + $innerXml + '</w:p>'
    $sub.InsertXML($xml)
}

Set-ParaRuns $d 2 '<w:proofErr w:type="spellStart"/><w:r><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a</w:t></w:r><w:r w:rsidRPr="000B6ACB"><w:t xml:space="preserve"> pure PHP library for reading and writing Word</w:t></w:r><w:r><w:t>, ODT, and RTF</w:t></w:r><w:r w:rsidRPr="000B6ACB"><w:t xml:space="preserve"> files</w:t></w:r><w:r><w:t xml:space="preserve">. This file is the </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:rStyle w:val="ChangedFontStyleChar"/></w:rPr><w:t>source file</w:t></w:r><w:r><w:t xml:space="preserve"> of read/write </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>capabilites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of PHP Word. The text in red should be changed when writing.</w:t></w:r>'
Set-ParaRuns $d 4 '<w:r><w:t>This should be deleted.</w:t></w:r>'
Set-ParaRuns $d 6 '<w:r><w:t>${CLONEME}</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-ParaRuns $d 7 '<w:proofErr w:type="spellStart"/><w:r><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can apply font formats such as </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:b/></w:rPr><w:t>bold</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:i/></w:rPr><w:t>italics</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>color</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:u w:val="dash"/></w:rPr><w:t>underline</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:strike/></w:rPr><w:t>strikethrough</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>subscript</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>superscript</w:t></w:r><w:r><w:t xml:space="preserve">, or </w:t></w:r><w:r w:rsidRPr="004348EF"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>highlighted</w:t></w:r><w:r><w:t xml:space="preserve">. You may also notice that there’s one text break (empty paragraph) before this one that can be created also by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>cloneValue</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t>}</w:t></w:r>'
Set-ParaRuns $d 9 '<w:proofErr w:type="spellStart"/><w:r><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can also format paragraph such as this justified, 12pt before and 12pt after with 1.5 lines spacing paragraph. This formatting can be applied inline or using predefined style as we use to do in Word.</w:t></w:r>'
Set-ParaRuns $d 11 '<w:proofErr w:type="spellStart"/><w:r><w:lastRenderedPageBreak/><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can also create multicolumn paragraph which is treated as a new section with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>continous</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> break type. We can define </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>either the number of columns or spacing between the columns.</w:t></w:r>'
Set-ParaRuns $d 13 '<w:proofErr w:type="spellStart"/><w:r><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can create bulleted lists …</w:t></w:r>'
Set-ParaRuns $d 16 '<w:r><w:t xml:space="preserve">… </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> numbered lists too.</w:t></w:r>'
Set-ParaRuns $d 20 '<w:r><w:t xml:space="preserve">Tables are also easy to be made by using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r>'
Set-ParaRuns $d 21 '<w:proofErr w:type="spellStart"/><w:r><w:t>PHPWord</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can insert images in your documents.</w:t></w:r>'

Write-Output 'DONE'
Write-Output $d.Paragraphs.Count
